# Atualização automática: 2025-08-30 21:00:25
# Applies the row-data refresh to the dashboard data sheet (rows 7-11 cyclic
# rotation of the per-detection fields, plus updated image/coord/confidence
# values for rows 16-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    # Force a Text number format for values that look numeric (e.g. "0.76")
    # so Excel keeps them stored as strings instead of converting to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
$ws.Range("D7").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E7").Value = "PLACA_20250717165933"
$ws.Range("F7").Value = "Beja"
$ws.Range("G7").Value = 38.02035
$ws.Range("H7").Value = -7.94715
Set-TextCell $ws "I7" "962,713,1006,765"
Set-TextCell $ws "J7" "0.76"

# --- Row 8 -----------------------------------------------------------------
$ws.Range("A8").Value = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
Set-TextCell $ws "I8" "967,614,1002,659"
Set-TextCell $ws "J8" "0.73"

# --- Row 9 -----------------------------------------------------------------
$ws.Range("A9").Value = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
Set-TextCell $ws "I9" "702,633,740,690"
Set-TextCell $ws "J9" "0.72"

# --- Row 10 ------------------------------------------------------------
$ws.Range("A10").Value = "dfd476d4-7689-4671-a076-78fe3ce806bb"
Set-TextCell $ws "I10" "1254,850,1294,895"
Set-TextCell $ws "J10" "0.67"

# --- Row 11 ------------------------------------------------------------
$ws.Range("A11").Value = "2117575c-4ae1-458c-b88a-fc40f40debdb"
$ws.Range("D11").Value = "image_20250727074723_ppp0.jpg"
$ws.Range("E11").Value = "PLACA_20250723145134"
$ws.Range("F11").Value = "Moura"
$ws.Range("G11").Value = 38.06587
$ws.Range("H11").Value = -7.221796
Set-TextCell $ws "I11" "1490,161,1563,258"
Set-TextCell $ws "J11" "0.62"

# --- Row 16 ------------------------------------------------------------
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
Set-TextCell $ws "I16" "643,531,686,575"
Set-TextCell $ws "J16" "0.76"

# --- Row 17 ------------------------------------------------------------
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
Set-TextCell $ws "I17" "794,481,830,526"
Set-TextCell $ws "J17" "0.72"

$wb.Save()
